$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-70, per diff (row 20 unchanged).
$gValues = @{
    2 = 3
    3 = 2
    4 = 3
    5 = 2
    6 = 1
    7 = 2
    8 = 1
    9 = 0
    10 = 1
    11 = 1
    12 = 2
    13 = 2
    14 = 2
    15 = 1
    16 = 2
    17 = 0
    18 = 3
    19 = 3
    21 = 2
    22 = 2
    23 = 2
    24 = 2
    25 = 3
    26 = 1
    27 = 1
    28 = 1
    29 = 1
    30 = 1
    31 = 0
    32 = 2
    33 = 0
    34 = 1
    35 = 1
    36 = 2
    37 = 1
    38 = 1
    39 = 3
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    44 = 1
    45 = 1
    46 = 4
    47 = 3
    48 = 1
    49 = 1
    50 = 1
    51 = 2
    52 = 5
    53 = 1
    54 = 1
    55 = 0
    56 = 1
    57 = 1
    58 = 1
    59 = 0
    60 = 2
    61 = 0
    62 = 1
    63 = 2
    64 = 2
    65 = 1
    66 = 3
    67 = 3
    68 = 3
    69 = 2
    70 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

$wb.Save()
